$d = $word.ActiveDocument

# The template used M2Doc "userdoc" fields (Word field codes) to mark an
# inlined user-content zone:
#   { m:userdoc 'zone1' }   ... protected content ...   { m:enduserdoc }
# The updated parser (TokenIteratorFieldRewriterSplit) expects these markers
# as plain literal text split across several runs instead of real Word
# fields, so every field is turned into the equivalent "{ ... }" text while
# preserving any bookmark that sits inside it (e.g. the automatic _GoBack
# bookmark Word drops into the last edited field).

# Fragment used to rebuild a paragraph's content as WordOpenXML for
# Range.InsertXML - $BODY$ is substituted with the paragraph-specific
# <w:r>/<w:bookmarkStart>/<w:bookmarkEnd> markup.
$xmlTemplate = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>$BODY$</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Runs (split the same way the original field's instrText runs were split)
# for the "{ m:userdoc 'zone1' }" field -> "{" / "m" / ":userdoc 'zone1'" / "}"
$zoneBody = "<w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc 'zone1'</w:t></w:r><w:r><w:t xml:space=`"preserve`">}</w:t></w:r>"

# Runs for the "{ m:enduserdoc }" field, keeping the _GoBack bookmark (it was
# sitting between the " m:" and "enduserdoc " instrText runs) in place ->
# "{m:" / <bookmark> / "enduserdoc}"
$endBody = "<w:r><w:t>{m:</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t xml:space=`"preserve`">enduserdoc}</w:t></w:r>"

# Phase 1: figure out, for every userdoc/enduserdoc field, which paragraph
# (1-based index) owns it and which replacement body it needs. This is done
# up front because deleting/replacing a field's host paragraph in place
# shrinks $d.Fields, which would otherwise shift the indices of the fields
# still to be processed. Paragraph indices, on the other hand, stay stable
# since InsertXML swaps a paragraph's content without adding/removing
# paragraphs. Field.Result/.Code "Paragraphs" chaining is not reliable in
# this host, so the owning paragraph is resolved via character offsets
# against the document-level collections instead.
$targets = @()
$fieldCount = $d.Fields.Count
for ($fi = 1; $fi -le $fieldCount; $fi++) {
    $f = $d.Fields.Item($fi)
    $code = $f.Code.Text
    $fStart = $f.Code.Start

    $body = $null
    if ($code -match "userdoc 'zone1'") {
        $body = $zoneBody
    } elseif ($code -match "enduserdoc") {
        $body = $endBody
    }

    if ($body -ne $null) {
        $paraCount = $d.Paragraphs.Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $p = $d.Paragraphs.Item($pi)
            if ($fStart -ge $p.Range.Start -and $fStart -lt $p.Range.End) {
                $targets += , @($pi, $body)
                break
            }
        }
    }
}

# Phase 2: apply the replacements by (stable) paragraph index.
foreach ($target in $targets) {
    $pi = $target[0]
    $body = $target[1]
    $p = $d.Paragraphs.Item($pi)
    $xml = $xmlTemplate -replace '\$BODY\$', $body
    $p.Range.InsertXML($xml)
}
